# Updates cryptos list values (price & 1h volume %) plus a row-order
# change for CoreDAO/Stellar/USDe/ThetaToken (rows 48-51), matching the
# "Updated cryptos list" GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "595.79", "66.441.43") that must
# stay plain text (matches the source inlineStr cells), so force Text format
# for the whole column before writing, then strip the formatting override
# back off afterwards so no stray styles are introduced.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value2 = "66.441.43"
$ws.Range("E2").Value2 = "  +2.22%  "
$ws.Range("D3").Value2 = "3.186.14"
$ws.Range("E3").Value2 = "  +0.67%  "
$ws.Range("E4").Value2 = "  +0.00%  "
$ws.Range("D5").Value2 = "595.79"
$ws.Range("E5").Value2 = "  +2.88%  "
$ws.Range("D6").Value2 = "154.45"
$ws.Range("E6").Value2 = "  +2.58%  "
$ws.Range("E7").Value2 = "  +0.06%  "
$ws.Range("D8").Value2 = "0.559"
$ws.Range("E8").Value2 = "  +5.68%  "
$ws.Range("D9").Value2 = "3.182.69"
$ws.Range("E9").Value2 = "  +0.55%  "
$ws.Range("E10").Value2 = "  +0.47%  "
$ws.Range("D11").Value2 = "5.85"
$ws.Range("E11").Value2 = "  -4.72%  "
$ws.Range("D12").Value2 = "0.518"
$ws.Range("E12").Value2 = "  +2.81%  "
$ws.Range("E13").Value2 = "  +0.39%  "
$ws.Range("D14").Value2 = "39.21"
$ws.Range("E14").Value2 = "  +4.45%  "
$ws.Range("D15").Value2 = "3.706.93"
$ws.Range("E15").Value2 = "  +0.62%  "
$ws.Range("D16").Value2 = "7.51"
$ws.Range("E16").Value2 = "  +4.42%  "
$ws.Range("D17").Value2 = "66.435.41"
$ws.Range("E17").Value2 = "  +2.25%  "
$ws.Range("D18").Value2 = "3.186.85"
$ws.Range("E18").Value2 = "  +1.33%  "
$ws.Range("E19").Value2 = "  +0.51%  "
$ws.Range("D20").Value2 = "519.23"
$ws.Range("E20").Value2 = "  +2.25%  "
$ws.Range("D21").Value2 = "15.43"
$ws.Range("E21").Value2 = "  +2.77%  "
$ws.Range("E22").Value2 = "  +2.51%  "
$ws.Range("D23").Value2 = "8.13"
$ws.Range("E23").Value2 = "  +4.53%  "
$ws.Range("D24").Value2 = "14.94"
$ws.Range("E24").Value2 = "  -2.53%  "
$ws.Range("D25").Value2 = "86.09"
$ws.Range("E25").Value2 = "  +1.64%  "
$ws.Range("E26").Value2 = "  -0.11%  "
$ws.Range("D27").Value2 = "9.26"
$ws.Range("E27").Value2 = "  +2.84%  "
$ws.Range("D28").Value2 = "3.00"
$ws.Range("E28").Value2 = "  +2.38%  "
$ws.Range("E29").Value2 = "  +6.41%  "
$ws.Range("D30").Value2 = "7.07"
$ws.Range("E30").Value2 = "  +12.22%  "
$ws.Range("E31").Value2 = "  +3.81%  "
$ws.Range("D32").Value2 = "28.32"
$ws.Range("E32").Value2 = "  +1.90%  "
$ws.Range("E33").Value2 = "  +2.36%  "
$ws.Range("E34").Value2 = "  +0.25%  "
$ws.Range("E35").Value2 = "  +0.01%  "
$ws.Range("D36").Value2 = "511.49"
$ws.Range("E36").Value2 = "  +6.31%  "
$ws.Range("D37").Value2 = "54.95"
$ws.Range("E37").Value2 = "  +0.12%  "
$ws.Range("D38").Value2 = "0.0902"
$ws.Range("E38").Value2 = "  +0.91%  "
$ws.Range("D39").Value2 = "0.0426"
$ws.Range("E39").Value2 = "  +1.64%  "
$ws.Range("D40").Value2 = "0.128"
$ws.Range("E40").Value2 = "  +10.12%  "
$ws.Range("D41").Value2 = "8.91"
$ws.Range("E41").Value2 = "  +1.14%  "
$ws.Range("E42").Value2 = "  -2.60%  "
$ws.Range("D43").Value2 = "0.301"
$ws.Range("E43").Value2 = "  +5.99%  "
$ws.Range("E44").Value2 = "  +12.60%  "
$ws.Range("E45").Value2 = "  +0.09%  "
$ws.Range("D46").Value2 = "2.901.55"
$ws.Range("E46").Value2 = "  -3.64%  "
$ws.Range("D47").Value2 = "28.52"
$ws.Range("E47").Value2 = "  -0.14%  "
$ws.Range("B48").Value2 = "CoreDAO"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D48").Value2 = "2.79"
$ws.Range("E48").Value2 = "  +11.09%  "
$ws.Range("B49").Value2 = "Stellar"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value2 = "0.118"
$ws.Range("E49").Value2 = "  +3.08%  "
$ws.Range("B50").Value2 = "USDe"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D50").Value2 = "0.999"
$ws.Range("E50").Value2 = "  -0.04%  "
$ws.Range("B51").Value2 = "ThetaToken"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value2 = "2.35"
$ws.Range("E51").Value2 = "  +4.25%  "

$ws.Range("D2:D51").ClearFormats()

